$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAllTypes   = $wb.Worksheets.Item("全部类型")

# sheet1 (展览)
$wsExhibition.Range("F4").Value = 114
$wsExhibition.Range("F5").Value = 1725
$wsExhibition.Range("F6").Value = 3303
$wsExhibition.Range("F7").Value = 978
$wsExhibition.Range("F8").Value = 2146
$wsExhibition.Range("F9").Value = 2065
$wsExhibition.Range("F10").Value = 1076
$wsExhibition.Range("F11").Value = 576
$wsExhibition.Range("F18").Value = 152
$wsExhibition.Range("F19").Value = 1522
$wsExhibition.Range("F21").Value = 677
$wsExhibition.Range("F23").Value = 12026
$wsExhibition.Range("F24").Value = 12036
$wsExhibition.Range("F29").Value = 301
$wsExhibition.Range("F30").Value = 1884
$wsExhibition.Range("F32").Value = 514

# sheet4 (全部类型)
$wsAllTypes.Range("F6").Value = 114
$wsAllTypes.Range("F7").Value = 1725
$wsAllTypes.Range("F8").Value = 3303
$wsAllTypes.Range("F9").Value = 978
$wsAllTypes.Range("F10").Value = 2146
$wsAllTypes.Range("F11").Value = 2065
$wsAllTypes.Range("F12").Value = 1076
$wsAllTypes.Range("F13").Value = 576
$wsAllTypes.Range("F22").Value = 152
$wsAllTypes.Range("F23").Value = 1522
$wsAllTypes.Range("F25").Value = 677
$wsAllTypes.Range("F27").Value = 12026
$wsAllTypes.Range("F28").Value = 12037
$wsAllTypes.Range("F33").Value = 301
$wsAllTypes.Range("F34").Value = 1884
$wsAllTypes.Range("F38").Value = 514
